$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("10000", "2.09074", "0.01683", "0.020335"),
    @("20000", "8.45289", "0.036186", "0.036298"),
    @("30000", "18.9352", "0.053639", "0.056458"),
    @("40000", "33.1753", "0.077115", "0.084508"),
    @("50000", "53.4929", "0.100999", "0.116662"),
    @("60000", "78.4559", "0.120201", "0.133587"),
    @("70000", "106.37", "0.14508", "0.179069"),
    @("80000", "137.933", "0.169194", "0.184616"),
    @("90000", "171.829", "0.191406", "0.231214"),
    @("100000", "217.711", "0.215667", "0.241235")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 1; $c -le 4; $c++) {
        $cell = $ws.Cells.Item($row, $c)
        $cell.Value = "'" + $values[$c - 1]
    }
}
